$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 44.09582266666666
$ws.Cells.Item(2,8).Value = 132.287468
$ws.Cells.Item(2,9).Value = 0.1927468402671175
$ws.Cells.Item(2,10).Value = 0.1927468402671175
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,13).Value = 5.802352333333334
$ws.Cells.Item(2,14).Value = 17.407057
$ws.Cells.Item(2,15).Value = 0.1988288201859468
$ws.Cells.Item(2,16).Value = 0.1988288201859468
$ws.Cells.Item(2,17).Value = 255.8594995401862
$ws.Cells.Item(2,18).Value = 2302.735495861676
$ws.Cells.Item(2,19).Value = 0.03832362684488012
$ws.Cells.Item(2,20).Value = 0.03832362684488012

$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 44.09582266666666
$ws.Cells.Item(3,8).Value = 132.287468
$ws.Cells.Item(3,9).Value = 0.1927468402671175
$ws.Cells.Item(3,10).Value = 0.1927468402671175
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 12.81923966666667
$ws.Cells.Item(3,14).Value = 38.457719
$ws.Cells.Item(3,15).Value = 0.4392760301648158
$ws.Cells.Item(3,16).Value = 0.4392760301648158
$ws.Cells.Item(3,17).Value = 565.2749190628324
$ws.Cells.Item(3,18).Value = 5087.474271565491
$ws.Cells.Item(3,19).Value = 0.08466906681935125
$ws.Cells.Item(3,20).Value = 0.08466906681935125

$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,7).Value = 44.09582266666666
$ws.Cells.Item(4,8).Value = 132.287468
$ws.Cells.Item(4,9).Value = 0.1927468402671175
$ws.Cells.Item(4,10).Value = 0.1927468402671175
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,13).Value = 3.182820333333333
$ws.Cells.Item(4,14).Value = 9.548461
$ws.Cells.Item(4,15).Value = 0.109065491956597
$ws.Cells.Item(4,16).Value = 0.109065491956597
$ws.Cells.Item(4,17).Value = 140.3490809985275
$ws.Cells.Item(4,18).Value = 1263.141728986748
$ws.Cells.Item(4,19).Value = 0.02102202895681279
$ws.Cells.Item(4,20).Value = 0.02102202895681279

$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,7).Value = 44.09582266666666
$ws.Cells.Item(5,8).Value = 132.287468
$ws.Cells.Item(5,9).Value = 0.1927468402671175
$ws.Cells.Item(5,10).Value = 0.1927468402671175
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,13).Value = 7.378240000000001
$ws.Cells.Item(5,14).Value = 22.13472
$ws.Cells.Item(5,15).Value = 0.2528296576926404
$ws.Cells.Item(5,16).Value = 0.2528296576926404
$ws.Cells.Item(5,17).Value = 325.3495626321067
$ws.Cells.Item(5,18).Value = 2928.14606368896
$ws.Cells.Item(5,19).Value = 0.04873211764607337
$ws.Cells.Item(5,20).Value = 0.04873211764607337

$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,7).Value = 111.3149186666667
$ws.Cells.Item(6,8).Value = 333.944756
$ws.Cells.Item(6,9).Value = 0.4865676055026886
$ws.Cells.Item(6,10).Value = 0.4865676055026886
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,13).Value = 5.802352333333334
$ws.Cells.Item(6,14).Value = 17.407057
$ws.Cells.Item(6,15).Value = 0.1988288201859468
$ws.Cells.Item(6,16).Value = 0.1988288201859468
$ws.Cells.Item(6,17).Value = 645.8883780603436
$ws.Cells.Item(6,18).Value = 5812.995402543092
$ws.Cells.Item(6,19).Value = 0.09674366294280076
$ws.Cells.Item(6,20).Value = 0.09674366294280076

$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,7).Value = 111.3149186666667
$ws.Cells.Item(7,8).Value = 333.944756
$ws.Cells.Item(7,9).Value = 0.4865676055026886
$ws.Cells.Item(7,10).Value = 0.4865676055026886
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,13).Value = 12.81923966666667
$ws.Cells.Item(7,14).Value = 38.457719
$ws.Cells.Item(7,15).Value = 0.4392760301648158
$ws.Cells.Item(7,16).Value = 0.4392760301648158
$ws.Cells.Item(7,17).Value = 1426.972620863507
$ws.Cells.Item(7,18).Value = 12842.75358777156
$ws.Cells.Item(7,19).Value = 0.2137374861520212
$ws.Cells.Item(7,20).Value = 0.2137374861520212

$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,7).Value = 111.3149186666667
$ws.Cells.Item(8,8).Value = 333.944756
$ws.Cells.Item(8,9).Value = 0.4865676055026886
$ws.Cells.Item(8,10).Value = 0.4865676055026886
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,13).Value = 3.182820333333333
$ws.Cells.Item(8,14).Value = 9.548461
$ws.Cells.Item(8,15).Value = 0.109065491956597
$ws.Cells.Item(8,16).Value = 0.109065491956597
$ws.Cells.Item(8,17).Value = 354.2953865356129
$ws.Cells.Item(8,18).Value = 3188.658478820516
$ws.Cells.Item(8,19).Value = 0.05306773526429414
$ws.Cells.Item(8,20).Value = 0.05306773526429414

$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,7).Value = 111.3149186666667
$ws.Cells.Item(9,8).Value = 333.944756
$ws.Cells.Item(9,9).Value = 0.4865676055026886
$ws.Cells.Item(9,10).Value = 0.4865676055026886
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,13).Value = 7.378240000000001
$ws.Cells.Item(9,14).Value = 22.13472
$ws.Cells.Item(9,15).Value = 0.2528296576926404
$ws.Cells.Item(9,16).Value = 0.2528296576926404
$ws.Cells.Item(9,17).Value = 821.3081855031467
$ws.Cells.Item(9,18).Value = 7391.77366952832
$ws.Cells.Item(9,19).Value = 0.1230187211435725
$ws.Cells.Item(9,20).Value = 0.1230187211435725

$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,7).Value = 40.34450033333334
$ws.Cells.Item(10,8).Value = 121.033501
$ws.Cells.Item(10,9).Value = 0.176349470111689
$ws.Cells.Item(10,10).Value = 0.176349470111689
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,13).Value = 5.802352333333334
$ws.Cells.Item(10,14).Value = 17.407057
$ws.Cells.Item(10,15).Value = 0.1988288201859468
$ws.Cells.Item(10,16).Value = 0.1988288201859468
$ws.Cells.Item(10,17).Value = 234.0930056462842
$ws.Cells.Item(10,18).Value = 2106.837050816557
$ws.Cells.Item(10,19).Value = 0.035063357082724
$ws.Cells.Item(10,20).Value = 0.035063357082724

$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,7).Value = 40.34450033333334
$ws.Cells.Item(11,8).Value = 121.033501
$ws.Cells.Item(11,9).Value = 0.176349470111689
$ws.Cells.Item(11,10).Value = 0.176349470111689
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,13).Value = 12.81923966666667
$ws.Cells.Item(11,14).Value = 38.457719
$ws.Cells.Item(11,15).Value = 0.4392760301648158
$ws.Cells.Item(11,16).Value = 0.4392760301648158
$ws.Cells.Item(11,17).Value = 517.1858190049132
$ws.Cells.Item(11,18).Value = 4654.672371044218
$ws.Cells.Item(11,19).Value = 0.07746609515233156
$ws.Cells.Item(11,20).Value = 0.07746609515233156

$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,7).Value = 40.34450033333334
$ws.Cells.Item(12,8).Value = 121.033501
$ws.Cells.Item(12,9).Value = 0.176349470111689
$ws.Cells.Item(12,10).Value = 0.176349470111689
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,13).Value = 3.182820333333333
$ws.Cells.Item(12,14).Value = 9.548461
$ws.Cells.Item(12,15).Value = 0.109065491956597
$ws.Cells.Item(12,16).Value = 0.109065491956597
$ws.Cells.Item(12,17).Value = 128.4092959991068
$ws.Cells.Item(12,18).Value = 1155.683663991961
$ws.Cells.Item(12,19).Value = 0.01923364171401656
$ws.Cells.Item(12,20).Value = 0.01923364171401656

$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,7).Value = 40.34450033333334
$ws.Cells.Item(13,8).Value = 121.033501
$ws.Cells.Item(13,9).Value = 0.176349470111689
$ws.Cells.Item(13,10).Value = 0.176349470111689
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,13).Value = 7.378240000000001
$ws.Cells.Item(13,14).Value = 22.13472
$ws.Cells.Item(13,15).Value = 0.2528296576926404
$ws.Cells.Item(13,16).Value = 0.2528296576926404
$ws.Cells.Item(13,17).Value = 297.6714061394134
$ws.Cells.Item(13,18).Value = 2679.04265525472
$ws.Cells.Item(13,19).Value = 0.04458637616261685
$ws.Cells.Item(13,20).Value = 0.04458637616261685

$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,7).Value = 33.02061066666666
$ws.Cells.Item(14,8).Value = 99.061832
$ws.Cells.Item(14,9).Value = 0.1443360841185049
$ws.Cells.Item(14,10).Value = 0.144336084118505
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,13).Value = 5.802352333333334
$ws.Cells.Item(14,14).Value = 17.407057
$ws.Cells.Item(14,15).Value = 0.1988288201859468
$ws.Cells.Item(14,16).Value = 0.1988288201859468
$ws.Cells.Item(14,17).Value = 191.5972173498249
$ws.Cells.Item(14,18).Value = 1724.374956148424
$ws.Cells.Item(14,19).Value = 0.02869817331554191
$ws.Cells.Item(14,20).Value = 0.02869817331554192

$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,7).Value = 33.02061066666666
$ws.Cells.Item(15,8).Value = 99.061832
$ws.Cells.Item(15,9).Value = 0.1443360841185049
$ws.Cells.Item(15,10).Value = 0.144336084118505
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,13).Value = 12.81923966666667
$ws.Cells.Item(15,14).Value = 38.457719
$ws.Cells.Item(15,15).Value = 0.4392760301648158
$ws.Cells.Item(15,16).Value = 0.4392760301648158
$ws.Cells.Item(15,17).Value = 423.2991220756897
$ws.Cells.Item(15,18).Value = 3809.692098681207
$ws.Cells.Item(15,19).Value = 0.06340338204111176
$ws.Cells.Item(15,20).Value = 0.06340338204111177

$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,7).Value = 33.02061066666666
$ws.Cells.Item(16,8).Value = 99.061832
$ws.Cells.Item(16,9).Value = 0.1443360841185049
$ws.Cells.Item(16,10).Value = 0.144336084118505
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,13).Value = 3.182820333333333
$ws.Cells.Item(16,14).Value = 9.548461
$ws.Cells.Item(16,15).Value = 0.109065491956597
$ws.Cells.Item(16,16).Value = 0.109065491956597
$ws.Cells.Item(16,17).Value = 105.0986710489502
$ws.Cells.Item(16,18).Value = 945.8880394405519
$ws.Cells.Item(16,19).Value = 0.01574208602147351
$ws.Cells.Item(16,20).Value = 0.01574208602147351

$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,7).Value = 33.02061066666666
$ws.Cells.Item(17,8).Value = 99.061832
$ws.Cells.Item(17,9).Value = 0.1443360841185049
$ws.Cells.Item(17,10).Value = 0.144336084118505
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,13).Value = 7.378240000000001
$ws.Cells.Item(17,14).Value = 22.13472
$ws.Cells.Item(17,15).Value = 0.2528296576926404
$ws.Cells.Item(17,16).Value = 0.2528296576926404
$ws.Cells.Item(17,17).Value = 243.6339904452267
$ws.Cells.Item(17,18).Value = 2192.70591400704
$ws.Cells.Item(17,19).Value = 0.03649244274037776
$ws.Cells.Item(17,20).Value = 0.03649244274037777

